# Delete row 11 ("Nos experimentos não ficou claro a divisão de treinamento e teste...")
# entirely, shifting all subsequent rows up by one. This matches the author's commit
# "erros a serem corrigidos" removing a resolved/duplicate entry from the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(11).Delete()

# Update the active view to match post-edit state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("A11").Select()
